# Parametrage_de_la_pelleteuse.xlsx - "Resolution pblm verin direction"
#
# Update the hard-coded vérin (cylinder) direction dimensions on Feuil1.
# Every other cell on the sheet (B75..B96 etc.) is a formula that
# recursively depends on these inputs, so changing the inputs is enough
# for Excel to recompute the whole cascade automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Main verin direction inputs -------------------------------------
$ws.Range("B89").Value = 1300
$ws.Range("B90").Value = 1000

# B93 used to be 1.5*B102 ; fix the ratio to 1.1*B102
$ws.Range("B93").Formula = "=1.1*B102"

# Another corrected dimension further down the sheet
$ws.Range("B106").Value = 280

# --- New parameter row -------------------------------------------------
$ws.Range("A111").Value = "Epaisseur_rotule_pivot_direction (mm)"
$ws.Range("B111").Value = 10

# --- Restore the on-screen selection/scroll position -------------------
$ws.Range("B93").Select()
$excel.ActiveWindow.ScrollRow = 87
$excel.ActiveWindow.ScrollColumn = 1
